$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 22
$ws1.Range("F6").Value = 15146
$ws1.Range("F7").Value = 410
$ws1.Range("F10").Value = 15274
$ws1.Range("F11").Value = 44
$ws1.Range("F12").Value = 8812
$ws1.Range("F13").Value = 348
$ws1.Range("F20").Value = 33
$ws1.Range("F21").Value = 527
$ws1.Range("F24").Value = 54
$ws1.Range("F27").Value = 15
$ws1.Range("F29").Value = 30
$ws1.Range("F32").Value = 35
$ws1.Range("F33").Value = 33
$ws1.Range("F34").Value = 234
$ws1.Range("F35").Value = 282
$ws1.Range("F36").Value = 432
$ws1.Range("F37").Value = 112
$ws1.Range("F38").Value = 5412

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1007

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 22
$ws4.Range("F6").Value = 15146
$ws4.Range("F7").Value = 410
$ws4.Range("F10").Value = 15274
$ws4.Range("F11").Value = 44
$ws4.Range("F12").Value = 8812
$ws4.Range("F13").Value = 348
$ws4.Range("F15").Value = 1007
$ws4.Range("F21").Value = 33
$ws4.Range("F22").Value = 527
$ws4.Range("F25").Value = 54
$ws4.Range("F28").Value = 15
$ws4.Range("F30").Value = 30
$ws4.Range("F35").Value = 35
$ws4.Range("F36").Value = 33
$ws4.Range("F37").Value = 234
$ws4.Range("F38").Value = 282
$ws4.Range("F39").Value = 432
$ws4.Range("F40").Value = 112
$ws4.Range("F41").Value = 5412
